# Apply the StructureDefinition-insight-id.xlsx update:
#  - bump Version to 6.0.0
#  - bump Date
#  - replace the "Contact" rows with "Publisher: Alvearie Team" and a new
#    "Jurisdiction: United States of America" row on the Metadata sheet
#  - update the Short/Definition for the root Extension row on the Elements sheet

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

$meta.Range("B9").Value = "Alvearie Team"

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# The old sheet had a duplicated "Contact" row at row 11; remove it so
# everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Insight Id"
$elements.Range("L2").Value = "ID for the insight as it is known by the source (or source system) that this insight was provided from"
